$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.4171829743026194
$ws.Cells.Item(2, 4).Value = 0.06224863125155711
$ws.Cells.Item(2, 5).Value = 0.1762322384878772
$ws.Cells.Item(2, 6).Value = 1.556187929967365
$ws.Cells.Item(2, 7).Value = 1.47106727210172
$ws.Cells.Item(2, 8).Value = 1.309920094614284
$ws.Cells.Item(2, 11).Value = 2.067073081219064
$ws.Cells.Item(2, 12).Value = 0.1537947137782396
$ws.Cells.Item(2, 14).Value = 1.360320642383932
$ws.Cells.Item(3, 3).Value = 0.4086188172729521
$ws.Cells.Item(3, 4).Value = 0.06300606459510227
$ws.Cells.Item(3, 5).Value = 0.1722182091423434
$ws.Cells.Item(3, 6).Value = 1.541267486656594
$ws.Cells.Item(3, 7).Value = 1.45604978341035
$ws.Cells.Item(3, 8).Value = 1.310024726463524
$ws.Cells.Item(3, 11).Value = 1.900597830129016
$ws.Cells.Item(3, 12).Value = 0.1499387334270637
$ws.Cells.Item(3, 14).Value = 1.378697834434309
$ws.Cells.Item(4, 3).Value = 0.4035999761287883
$ws.Cells.Item(4, 4).Value = 0.06349179032226715
$ws.Cells.Item(4, 5).Value = 0.1698578013322987
$ws.Cells.Item(4, 6).Value = 1.533228010338803
$ws.Cells.Item(4, 7).Value = 1.447969616184508
$ws.Cells.Item(4, 8).Value = 1.310878517365339
$ws.Cells.Item(4, 11).Value = 1.799132840006735
$ws.Cells.Item(4, 12).Value = 0.1476645686336866
$ws.Cells.Item(4, 14).Value = 1.390573083689631
$ws.Cells.Item(5, 3).Value = 0.4016148542617941
$ws.Cells.Item(5, 4).Value = 0.06369491751197831
$ws.Cells.Item(5, 5).Value = 0.1689220397725038
$ws.Cells.Item(5, 6).Value = 1.530232743457134
$ws.Cells.Item(5, 7).Value = 1.444962250334427
$ws.Cells.Item(5, 8).Value = 1.311424461118648
$ws.Cells.Item(5, 11).Value = 1.757973639566728
$ws.Cells.Item(5, 12).Value = 0.1467612277279216
$ws.Cells.Item(5, 14).Value = 1.395560991843297
$ws.Cells.Item(6, 3).Value = 0.401288851694801
$ws.Cells.Item(6, 4).Value = 0.06372895995871186
$ws.Cells.Item(6, 5).Value = 0.1687682326762463
$ws.Cells.Item(6, 6).Value = 1.529752311999459
$ws.Cells.Item(6, 7).Value = 1.444480072762389
$ws.Cells.Item(6, 8).Value = 1.311527059295287
$ws.Cells.Item(6, 11).Value = 1.751150567484103
$ws.Cells.Item(6, 12).Value = 0.1466126392810665
$ws.Cells.Item(6, 14).Value = 1.396398205078963
$ws.Cells.Item(7, 3).Value = 0.4035729609454108
$ws.Cells.Item(7, 4).Value = 0.06349450875994744
$ws.Cells.Item(7, 5).Value = 0.1698450756592109
$ws.Cells.Item(7, 6).Value = 1.533186479405757
$ws.Cells.Item(7, 7).Value = 1.447927904282011
$ws.Cells.Item(7, 8).Value = 1.310885079091463
$ws.Cells.Item(7, 11).Value = 1.798576989427943
$ws.Cells.Item(7, 12).Value = 0.1476522912375415
$ws.Cells.Item(7, 14).Value = 1.390639750647644
$ws.Cells.Item(8, 3).Value = 0.4141802325071353
$ws.Cells.Item(8, 4).Value = 0.06250550370366614
$ws.Cells.Item(8, 5).Value = 0.1748265172822201
$ws.Cells.Item(8, 6).Value = 1.550809819811548
$ws.Cells.Item(8, 7).Value = 1.465651617417393
$ws.Cells.Item(8, 8).Value = 1.309791982407972
$ws.Cells.Item(8, 11).Value = 2.00951616187632
$ws.Cells.Item(8, 12).Value = 0.15244572442964
$ws.Cells.Item(8, 14).Value = 1.366534249905826
$ws.Cells.Item(9, 3).Value = 0.4368905296990988
$ws.Cells.Item(9, 4).Value = 0.06073017618086851
$ws.Cells.Item(9, 5).Value = 0.1854266279350156
$ws.Cells.Item(9, 6).Value = 1.594327423827593
$ws.Cells.Item(9, 7).Value = 1.509527994193462
$ws.Cells.Item(9, 8).Value = 1.313938548016523
$ws.Cells.Item(9, 11).Value = 2.429170573197553
$ws.Cells.Item(9, 12).Value = 0.1625917822963885
$ws.Cells.Item(9, 14).Value = 1.323961849961119
$ws.Cells.Item(10, 3).Value = 0.4547540756375099
$ws.Cells.Item(10, 4).Value = 0.05952612110751243
$ws.Cells.Item(10, 5).Value = 0.1937288889394821
$ws.Cells.Item(10, 6).Value = 1.631847389265786
$ws.Cells.Item(10, 7).Value = 1.547426771460152
$ws.Cells.Item(10, 8).Value = 1.320857666539695
$ws.Cells.Item(10, 11).Value = 2.741241878747985
$ws.Cells.Item(10, 12).Value = 0.1705087709406996
$ws.Cells.Item(10, 14).Value = 1.295553775690117
$ws.Cells.Item(11, 3).Value = 0.463139798675769
$ws.Cells.Item(11, 4).Value = 0.05900019708403104
$ws.Cells.Item(11, 5).Value = 0.1976191909292382
$ws.Cells.Item(11, 6).Value = 1.650140004765987
$ws.Cells.Item(11, 7).Value = 1.565920295897911
$ws.Cells.Item(11, 8).Value = 1.324854762089529
$ws.Cells.Item(11, 11).Value = 2.884048136147157
$ws.Cells.Item(11, 12).Value = 0.1742126523890875
$ws.Cells.Item(11, 14).Value = 1.283254713982021
$ws.Cells.Item(12, 3).Value = 0.4663527937632068
$ws.Cells.Item(12, 4).Value = 0.058804186420093
$ws.Cells.Item(12, 5).Value = 0.1991087957738245
$ws.Cells.Item(12, 6).Value = 1.657244503219886
$ws.Cells.Item(12, 7).Value = 1.573105277711278
$ws.Cells.Item(12, 8).Value = 1.32649118418496
$ws.Cells.Item(12, 11).Value = 2.938247576645949
$ws.Cells.Item(12, 12).Value = 0.1756300673031177
$ws.Cells.Item(12, 14).Value = 1.278687260309916
$ws.Cells.Item(13, 3).Value = 0.4656591467691555
$ws.Cells.Item(13, 4).Value = 0.05884626077567834
$ws.Cells.Item(13, 5).Value = 0.1987872502609704
$ws.Cells.Item(13, 6).Value = 1.655706506545798
$ws.Cells.Item(13, 7).Value = 1.571549746977297
$ws.Cells.Item(13, 8).Value = 1.326133279094222
$ws.Cells.Item(13, 11).Value = 2.92656933107628
$ws.Cells.Item(13, 12).Value = 0.1753241400621022
$ws.Cells.Item(13, 14).Value = 1.279666939792939
$ws.Cells.Item(14, 3).Value = 0.4634033809969651
$ws.Cells.Item(14, 4).Value = 0.05898400807458959
$ws.Cells.Item(14, 5).Value = 0.1977414117054224
$ws.Cells.Item(14, 6).Value = 1.65072093112677
$ws.Cells.Item(14, 7).Value = 1.566507753736346
$ws.Cells.Item(14, 8).Value = 1.324986926527913
$ws.Cells.Item(14, 11).Value = 2.888504712338772
$ws.Cells.Item(14, 12).Value = 0.1743289660068257
$ws.Cells.Item(14, 14).Value = 1.28287714313603
$ws.Cells.Item(15, 3).Value = 0.462026549050762
$ws.Cells.Item(15, 4).Value = 0.05906879218324246
$ws.Cells.Item(15, 5).Value = 0.1971029483703148
$ws.Cells.Item(15, 6).Value = 1.647690277048071
$ws.Cells.Item(15, 7).Value = 1.563443124302012
$ws.Cells.Item(15, 8).Value = 1.324300765079016
$ws.Cells.Item(15, 11).Value = 2.865204913761488
$ws.Cells.Item(15, 12).Value = 0.1737213284645094
$ws.Cells.Item(15, 14).Value = 1.284855204021014
$ws.Cells.Item(16, 3).Value = 0.4542112867400192
$ws.Cells.Item(16, 4).Value = 0.0595609310911831
$ws.Cells.Item(16, 5).Value = 0.1934769423707152
$ws.Cells.Item(16, 6).Value = 1.630676683975636
$ws.Cells.Item(16, 7).Value = 1.54624353651667
$ws.Cells.Item(16, 8).Value = 1.320613596969196
$ws.Cells.Item(16, 11).Value = 2.731926170144391
$ws.Cells.Item(16, 12).Value = 0.1702687835484795
$ws.Cells.Item(16, 14).Value = 1.296370121284397
$ws.Cells.Item(17, 3).Value = 0.4494834740010276
$ws.Cells.Item(17, 4).Value = 0.05986843422770427
$ws.Cells.Item(17, 5).Value = 0.1912816628789002
$ws.Cells.Item(17, 6).Value = 1.620554025417206
$ws.Cells.Item(17, 7).Value = 1.536014305585866
$ws.Cells.Item(17, 8).Value = 1.318569685038938
$ws.Cells.Item(17, 11).Value = 2.650380365352078
$ws.Cells.Item(17, 12).Value = 0.1681770646540031
$ws.Cells.Item(17, 14).Value = 1.303594099988089
$ws.Cells.Item(18, 3).Value = 0.4467885756853036
$ws.Cells.Item(18, 4).Value = 0.06004735384142101
$ws.Cells.Item(18, 5).Value = 0.1900296789416203
$ws.Cells.Item(18, 6).Value = 1.614846920772095
$ws.Cells.Item(18, 7).Value = 1.530248576324652
$ws.Cells.Item(18, 8).Value = 1.317474005924822
$ws.Cells.Item(18, 11).Value = 2.603556731417484
$ws.Cells.Item(18, 12).Value = 0.1669835952932743
$ws.Cells.Item(18, 14).Value = 1.307807853687642
$ws.Cells.Item(19, 3).Value = 0.4458803170081467
$ws.Cells.Item(19, 4).Value = 0.06010828508735511
$ws.Cells.Item(19, 5).Value = 0.189607610396628
$ws.Cells.Item(19, 6).Value = 1.612934333614092
$ws.Cells.Item(19, 7).Value = 1.528316590069949
$ws.Cells.Item(19, 8).Value = 1.317116735515214
$ws.Cells.Item(19, 11).Value = 2.587716667779432
$ws.Cells.Item(19, 12).Value = 0.1665811580680412
$ws.Cells.Item(19, 14).Value = 1.309244635731449
$ws.Cells.Item(20, 3).Value = 0.449984229584544
$ws.Cells.Item(20, 4).Value = 0.05983548756084556
$ws.Cells.Item(20, 5).Value = 0.1915142477703498
$ws.Cells.Item(20, 6).Value = 1.621619668964883
$ws.Cells.Item(20, 7).Value = 1.537091015179186
$ws.Cells.Item(20, 8).Value = 1.318778986422672
$ws.Cells.Item(20, 11).Value = 2.659052838765831
$ws.Cells.Item(20, 12).Value = 0.1683987340316691
$ws.Cells.Item(20, 14).Value = 1.302819017125909
$ws.Cells.Item(21, 3).Value = 0.4640649347313683
$ws.Cells.Item(21, 4).Value = 0.05894346291021524
$ws.Cells.Item(21, 5).Value = 0.1980481530614995
$ws.Cells.Item(21, 6).Value = 1.652180487537194
$ws.Cells.Item(21, 7).Value = 1.567983759576265
$ws.Cells.Item(21, 8).Value = 1.325320299889881
$ws.Cells.Item(21, 11).Value = 2.899681903298813
$ws.Cells.Item(21, 12).Value = 0.174620869236449
$ws.Cells.Item(21, 14).Value = 1.281931786042993
$ws.Cells.Item(22, 3).Value = 0.4734861800993428
$ws.Cells.Item(22, 4).Value = 0.05837881156011449
$ws.Cells.Item(22, 5).Value = 0.2024142643726918
$ws.Cells.Item(22, 6).Value = 1.673188905027502
$ws.Cells.Item(22, 7).Value = 1.589234913915874
$ws.Cells.Item(22, 8).Value = 1.330311555514584
$ws.Cells.Item(22, 11).Value = 3.057657182368132
$ws.Cells.Item(22, 12).Value = 0.1787739187314799
$ws.Cells.Item(22, 14).Value = 1.268805084739512
$ws.Cells.Item(23, 3).Value = 0.4684378109948284
$ws.Cells.Item(23, 4).Value = 0.05867849535286762
$ws.Cells.Item(23, 5).Value = 0.2000751862760737
$ws.Cells.Item(23, 6).Value = 1.661881137967498
$ws.Cells.Item(23, 7).Value = 1.577795135327534
$ws.Cells.Item(23, 8).Value = 1.327581887307019
$ws.Cells.Item(23, 11).Value = 2.973277629076563
$ws.Cells.Item(23, 12).Value = 0.1765494042395375
$ws.Cells.Item(23, 14).Value = 1.275763003836644
$ws.Cells.Item(24, 3).Value = 0.449757765786245
$ws.Cells.Item(24, 4).Value = 0.05985037610735056
$ws.Cells.Item(24, 5).Value = 0.1914090646542803
$ws.Cells.Item(24, 6).Value = 1.62113754104692
$ws.Cells.Item(24, 7).Value = 1.536603876086815
$ws.Cells.Item(24, 8).Value = 1.31868411400626
$ws.Cells.Item(24, 11).Value = 2.655131832495897
$ws.Cells.Item(24, 12).Value = 0.1682984890088761
$ws.Cells.Item(24, 14).Value = 1.303169243358781
$ws.Cells.Item(25, 3).Value = 0.4305409150396429
$ws.Cells.Item(25, 4).Value = 0.06119285818536291
$ws.Cells.Item(25, 5).Value = 0.1824692348788233
$ws.Cells.Item(25, 6).Value = 1.581586840172463
$ws.Cells.Item(25, 7).Value = 1.496671316132165
$ws.Cells.Item(25, 8).Value = 1.312139495213188
$ws.Cells.Item(25, 11).Value = 2.314990935877461
$ws.Cells.Item(25, 12).Value = 0.1597663240572871
$ws.Cells.Item(25, 14).Value = 1.334975273053079
